$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245 (shifts existing rows 245-387 down to 246-388,
# expanding the used range from A1:R387 to A1:R388).
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new weekly data point.
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 44830
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = 100112044
$ws.Range("G245").Value = "Perejil"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 55
$ws.Range("K245").Value = 3300
$ws.Range("L245").Value = 3300
$ws.Range("M245").Value = 3300
$ws.Range("N245").Value = '$/docena de atados (3 kilos)'
$ws.Range("O245").Value = "Región Metropolitana"
$ws.Range("P245").Value = 1100
$ws.Range("Q245").Value = 3
$ws.Range("R245").Value = "Hortaliza"
